# Generate Report for Handback
# The be69c26c-4b22-4ac7-93db-fdef1d2accaf item has come back from handback
# processing "in sync with en-US" - update its status everywhere, and record
# the handback target file + datetime on the per-language detail sheets.
# The d32f0acf-ec09-4e00-94d2-99390416cc32 item is still awaiting handoff,
# so it keeps its "Ready for handoff" status (unchanged).

$wb = $excel.ActiveWorkbook

$handedBackStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: file -> zh-cn status / de-de status / latest handoff date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 2 = be69c26c-4b22-4ac7-93db-fdef1d2accaf.md -> now handed back in both languages
$wsOverview.Range("B2").Value = $handedBackStatus
$wsOverview.Range("C2").Value = $handedBackStatus
# D2 (latest handoff date) is unchanged

# Row 3 = d32f0acf-ec09-4e00-94d2-99390416cc32.md -> still ready for handoff (no change)

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 2: be69c26c item - update status, add target/handback file links, stamp handback datetime
$wsZhCn.Range("C2").Value = $handedBackStatus

$wsZhCn.Range("F2").Value = "be69c26c-4b22-4ac7-93db-fdef1d2accaf.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/160c89485cb434c67ae3cd7e772399d89dac26ff/e2e/be69c26c-4b22-4ac7-93db-fdef1d2accaf.md", "", "", "be69c26c-4b22-4ac7-93db-fdef1d2accaf.md")

$wsZhCn.Range("G2").Value = "be69c26c-4b22-4ac7-93db-fdef1d2accaf.405238ea83a6a213a5847ce8d6b0ae90920e00ef.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c0dec2a97277437524baf296856b831f1fe12e50/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/be69c26c-4b22-4ac7-93db-fdef1d2accaf.405238ea83a6a213a5847ce8d6b0ae90920e00ef.zh-cn.xlf", "", "", "be69c26c-4b22-4ac7-93db-fdef1d2accaf.405238ea83a6a213a5847ce8d6b0ae90920e00ef.zh-cn.xlf")

$wsZhCn.Range("H2").Value = "2016-03-25 00:49:36"

# Row 3: d32f0acf item - still ready for handoff (no change)

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 2: be69c26c item - update status, add target/handback file links, stamp handback datetime
$wsDeDe.Range("C2").Value = $handedBackStatus

$wsDeDe.Range("F2").Value = "be69c26c-4b22-4ac7-93db-fdef1d2accaf.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/160c89485cb434c67ae3cd7e772399d89dac26ff/e2e/be69c26c-4b22-4ac7-93db-fdef1d2accaf.md", "", "", "be69c26c-4b22-4ac7-93db-fdef1d2accaf.md")

$wsDeDe.Range("G2").Value = "be69c26c-4b22-4ac7-93db-fdef1d2accaf.405238ea83a6a213a5847ce8d6b0ae90920e00ef.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/935e4ccb2a9f73dbbcfec046cdfd60b0cca60940/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/be69c26c-4b22-4ac7-93db-fdef1d2accaf.405238ea83a6a213a5847ce8d6b0ae90920e00ef.de-de.xlf", "", "", "be69c26c-4b22-4ac7-93db-fdef1d2accaf.405238ea83a6a213a5847ce8d6b0ae90920e00ef.de-de.xlf")

$wsDeDe.Range("H2").Value = "2016-03-25 00:49:44"

# Row 3: d32f0acf item - still ready for handoff (no change)

Write-Host "Handback report updated."
